$d = $word.ActiveDocument

# Locate the answer paragraphs that follow "11. What are the various
# components in your first app that you designed?" - the three trailing
# paragraphs holding the old (three-sentence) answer text.
$p70 = $d.Paragraphs.Item(70)
$p71 = $d.Paragraphs.Item(71)

# Drop the last two answer paragraphs entirely ("When the score..." and
# "Obstacle come...") - their content goes away, merging everything into
# a single remaining answer paragraph.
$delRange = $d.Range($p70.Range.Start, $p71.Range.End)
$delRange.Delete()

$p69 = $d.Paragraphs.Item(69)

# Mark the paragraph (and its run) as explicitly not-underlined, using the
# full paragraph range so both the paragraph mark's rPr and the run's rPr
# pick up <w:u w:val="none"/>.
$p69.Range.Font.Underline = 0

# Replace the old answer text with the new answer text, appending a
# temporary sentinel character so we can park the bookmark exactly at the
# end of the real text without landing on the paragraph-mark boundary
# (the COM host mis-resolves a bookmark collapsed exactly at end-of-run).
$r = $p69.Range.Duplicate
[void]$r.MoveEnd(1, -1)
$r.Text = "I learned how to make button. How to add alert for conformation, add text to the button, change colour of the button.X"

# Select just the sentinel "X" and collapse to its start - that is exactly
# the insertion point right after the real text.
$sentinel = $p69.Range.Duplicate
[void]$sentinel.MoveEnd(1, -1)
[void]$sentinel.MoveStart(1, $sentinel.End - $sentinel.Start - 1)
$sentinel.Collapse(1)

# The _GoBack bookmark used to sit in the empty paragraph right after the
# question; move it to the end of the (new) answer text.
$d.Bookmarks.Add("_GoBack", $sentinel)

# Remove the sentinel "X" now that the bookmark is anchored.
$sentinel2 = $p69.Range.Duplicate
[void]$sentinel2.MoveEnd(1, -1)
[void]$sentinel2.MoveStart(1, $sentinel2.End - $sentinel2.Start - 1)
$sentinel2.Delete()
